{"js": "// Replace the date line and the 25 division problems in the table with\n// the new values from the updated worksheet. Every old string below is\n// unique within the document, so a simple search-and-replace for each\n// pair is safe and keeps the original run formatting (search results\n// are themselves Range objects, and insertText(..., \"Replace\") swaps\n// only the text of the matched range, leaving rPr/pPr untouched).\nconst replacements = [\n  [\"2024-06-28 Friday\", \"2024-06-29 Saturday\"],\n  [\"484\u00f76=\", \"111\u00f74=\"],\n  [\"120\u00f79=\", \"192\u00f74=\"],\n  [\"930\u00f76=\", \"750\u00f77=\"],\n  [\"712\u00f73=\", \"226\u00f77=\"],\n  [\"421\u00f79=\", \"405\u00f77=\"],\n  [\"113\u00f74=\", \"907\u00f79=\"],\n  [\"466\u00f79=\", \"845\u00f79=\"],\n  [\"855\u00f77=\", \"981\u00f77=\"],\n  [\"847\u00f79=\", \"423\u00f72=\"],\n  [\"726\u00f75=\", \"476\u00f75=\"],\n  [\"855\u00f79=\", \"594\u00f77=\"],\n  [\"300\u00f72=\", \"408\u00f79=\"],\n  [\"154\u00f74=\", \"199\u00f78=\"],\n  [\"750\u00f76=\", \"621\u00f76=\"],\n  [\"397\u00f75=\", \"822\u00f73=\"],\n  [\"829\u00f75=\", \"344\u00f72=\"],\n  [\"474\u00f78=\", \"336\u00f77=\"],\n  [\"462\u00f79=\", \"665\u00f72=\"],\n  [\"235\u00f76=\", \"286\u00f75=\"],\n  [\"800\u00f75=\", \"496\u00f76=\"],\n  [\"868\u00f73=\", \"240\u00f74=\"],\n  [\"487\u00f79=\", \"139\u00f77=\"],\n  [\"187\u00f77=\", \"825\u00f72=\"],\n  [\"429\u00f78=\", \"959\u00f72=\"],\n  [\"471\u00f79=\", \"706\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 division problems in the table with\n# the new values from the updated worksheet. Word's Find/Execute with\n# Replace:=wdReplaceAll swaps only the matched text, leaving the run's\n# formatting (rPr/pPr) untouched, and every \"old\" string here is unique\n# within the document so there is no risk of an unintended match.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-06-28 Friday\", \"2024-06-29 Saturday\"),\n    @(\"484\u00f76=\", \"111\u00f74=\"),\n    @(\"120\u00f79=\", \"192\u00f74=\"),\n    @(\"930\u00f76=\", \"750\u00f77=\"),\n    @(\"712\u00f73=\", \"226\u00f77=\"),\n    @(\"421\u00f79=\", \"405\u00f77=\"),\n    @(\"113\u00f74=\", \"907\u00f79=\"),\n    @(\"466\u00f79=\", \"845\u00f79=\"),\n    @(\"855\u00f77=\", \"981\u00f77=\"),\n    @(\"847\u00f79=\", \"423\u00f72=\"),\n    @(\"726\u00f75=\", \"476\u00f75=\"),\n    @(\"855\u00f79=\", \"594\u00f77=\"),\n    @(\"300\u00f72=\", \"408\u00f79=\"),\n    @(\"154\u00f74=\", \"199\u00f78=\"),\n    @(\"750\u00f76=\", \"621\u00f76=\"),\n    @(\"397\u00f75=\", \"822\u00f73=\"),\n    @(\"829\u00f75=\", \"344\u00f72=\"),\n    @(\"474\u00f78=\", \"336\u00f77=\"),\n    @(\"462\u00f79=\", \"665\u00f72=\"),\n    @(\"235\u00f76=\", \"286\u00f75=\"),\n    @(\"800\u00f75=\", \"496\u00f76=\"),\n    @(\"868\u00f73=\", \"240\u00f74=\"),\n    @(\"487\u00f79=\", \"139\u00f77=\"),\n    @(\"187\u00f77=\", \"825\u00f72=\"),\n    @(\"429\u00f78=\", \"959\u00f72=\"),\n    @(\"471\u00f79=\", \"706\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
